$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "2019年5月14日23:24:48"
$ws.Range("B50").Value = "周二"
$ws.Range("C50").Value = "web的搭建"
$ws.Range("D50").Value = "22:00--01:30"

$ws.Range("A51").Value = "2019年5月15日22:39:24"
$ws.Range("B51").Value = "周三"
$ws.Range("C51").Value = "OSALS项目架构的搭建，web的完善，成功"
$ws.Range("D51").Value = "12:30--15:00 & 17:00--21:00"

$ws.Range("D52").Select()
